$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("J2").Value = 4.073E-02
$ws.Range("K2").Value = 2.525E-02
